$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

$ws.Range("D2").Value = "42.992.04"
$ws.Range("E2").Value = "  -0.99%  "
$ws.Range("D3").Value = "2.243.10"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.51%  "
Set-TextValue "D5" "113.92"
$ws.Range("E5").Value = "  +0.44%  "
Set-TextValue "D6" "273.99"
$ws.Range("E6").Value = "  +3.21%  "
Set-TextValue "D7" "0.627"
$ws.Range("E7").Value = "  +0.87%  "
Set-TextValue "D8" "1.01"
$ws.Range("E8").Value = "  +0.32%  "
Set-TextValue "D9" "0.608"
$ws.Range("E9").Value = "  +0.21%  "
Set-TextValue "D10" "46.35"
$ws.Range("E10").Value = "  -2.03%  "
Set-TextValue "D11" "0.0931"
$ws.Range("E11").Value = "  -0.12%  "
Set-TextValue "D12" "9.16"
$ws.Range("E12").Value = "  -1.35%  "
$ws.Range("E13").Value = "  -2.52%  "
Set-TextValue "D14" "15.37"
$ws.Range("E14").Value = "  -0.70%  "
Set-TextValue "D15" "0.875"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").Value = "2.582.01"
$ws.Range("E16").Value = "  -1.62%  "
$ws.Range("D17").Value = "2.249.88"
$ws.Range("E17").Value = "  -1.09%  "
$ws.Range("D18").Value = "43.008.63"
$ws.Range("E18").Value = "  -0.62%  "
$ws.Range("E19").Value = "  -1.12%  "
Set-TextValue "D20" "6.77"
$ws.Range("E20").Value = "  -0.07%  "
Set-TextValue "D21" "72.14"
$ws.Range("E21").Value = "  +0.10%  "
$ws.Range("E22").Value = "  -5.11%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue "D23" "2.98"
$ws.Range("E23").Value = "  +4.46%  "
$ws.Range("B24").Value = "BitcoinCash"
$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D24" "231.44"
$ws.Range("E24").Value = "  -0.99%  "
Set-TextValue "D25" "9.34"
$ws.Range("E25").Value = "  -2.58%  "
Set-TextValue "D26" "12.15"
$ws.Range("E26").Value = "  +6.12%  "
$ws.Range("E27").Value = "  -0.69%  "
Set-TextValue "D28" "40.57"
$ws.Range("E28").Value = "  -1.03%  "
Set-TextValue "D29" "2.25"
$ws.Range("E29").Value = "  -0.34%  "
Set-TextValue "D30" "3.27"
$ws.Range("E30").Value = "  -2.34%  "
Set-TextValue "D31" "174.08"
$ws.Range("E31").Value = "  +0.20%  "
Set-TextValue "D32" "21.16"
$ws.Range("E32").Value = "  -2.15%  "
Set-TextValue "D33" "0.0907"
$ws.Range("E33").Value = "  +0.26%  "
Set-TextValue "D34" "5.60"
$ws.Range("E34").Value = "  -2.51%  "
Set-TextValue "D35" "4.36"
$ws.Range("E35").Value = "  +10.99%  "
Set-TextValue "D36" "0.128"
$ws.Range("E36").Value = "  +0.00%  "
Set-TextValue "D37" "4.71"
$ws.Range("E37").Value = "  +1.66%  "
Set-TextValue "D38" "0.0374"
$ws.Range("E38").Value = "  +1.62%  "
Set-TextValue "D39" "0.107"
$ws.Range("E39").Value = "  +2.96%  "
Set-TextValue "D40" "2.57"
$ws.Range("E40").Value = "  -1.60%  "
Set-TextValue "D41" "71.22"
$ws.Range("E41").Value = "  -6.19%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue "D42" "0.233"
$ws.Range("E42").Value = "  -2.36%  "
$ws.Range("B43").Value = "Celestia"
$ws.Range("C43").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextValue "D43" "13.19"
$ws.Range("E43").Value = "  -6.48%  "
Set-TextValue "D44" "1.00"
$ws.Range("E44").Value = "  +0.05%  "
Set-TextValue "D45" "5.67"
$ws.Range("E45").Value = "  -7.87%  "
Set-TextValue "D46" "1.34"
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("B47").Value = "TrustWalletToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D47" "1.26"
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
Set-TextValue "D48" "8.44"
$ws.Range("E48").Value = "  -1.72%  "
$ws.Range("E49").Value = "  -1.00%  "
$ws.Range("B50").Value = "TheSandbox"
$ws.Range("C50").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
Set-TextValue "D50" "0.646"
$ws.Range("E50").Value = "  +8.00%  "
$ws.Range("B51").Value = "Aave"
$ws.Range("C51").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue "D51" "100.61"
$ws.Range("E51").Value = "  -2.80%  "
